$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 61
$ws.Range("N2").Value = 85.8724807945396

$ws.Range("K3").Value = 57.2
$ws.Range("N3").Value = 85.8724807945396

$ws.Range("K4").Value = 55.8
$ws.Range("N4").Value = 85.8724807945396

$ws.Range("K5").Value = 51.2
$ws.Range("N5").Value = 85.8724807945396

$ws.Range("K6").Value = 51
$ws.Range("N6").Value = 85.8724807945396

$ws.Range("K7").Value = 44.8
$ws.Range("N7").Value = 85.8724807945396
